$wb = $excel.ActiveWorkbook

# Rename the "Include from HLA HGNC GeneID " sheet to "Include #0"
$wsInclude = $wb.Worksheets.Item("Include from HLA HGNC GeneID ")
$wsInclude.Name = "Include #0"

# Metadata sheet updates
$wsMeta = $wb.Worksheets.Item("Metadata")

# Update Version value
$wsMeta.Range("B3").Value = "0.1.1"

# Update Date value
$wsMeta.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# Insert a new row at position 11 for "Jurisdiction" (shifts Description..Immutable down by one)
$wsMeta.Rows.Item(11).Insert()

$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""
